$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "81.449.77"
$ws.Range("E2").Value = "  +3.39%  "
$ws.Range("D3").Value = "3.180.87"
$ws.Range("E3").Value = "  +0.15%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "'208.00"
$ws.Range("E5").Value = "  +1.65%  "
$ws.Range("D6").Value = "'633.39"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("E7").Value = "  +28.46%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  +2.60%  "
$ws.Range("D10").Value = "3.177.41"
$ws.Range("E10").Value = "  +0.38%  "
$ws.Range("D11").Value = "'0.595"
$ws.Range("E11").Value = "  +4.91%  "
$ws.Range("D12").Value = "'0.0000262"
$ws.Range("E12").Value = "  +14.19%  "
$ws.Range("E13").Value = "  +1.97%  "
$ws.Range("E14").Value = "  -1.44%  "
$ws.Range("D15").Value = "3.763.40"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("D16").Value = "'31.93"
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("D17").Value = "81.478.45"
$ws.Range("E17").Value = "  +4.01%  "
$ws.Range("D18").Value = "3.179.63"
$ws.Range("E18").Value = "  +0.66%  "
$ws.Range("D19").Value = "'3.21"
$ws.Range("E19").Value = "  +12.28%  "
$ws.Range("D20").Value = "'14.17"
$ws.Range("E20").Value = "  -1.52%  "
$ws.Range("D21").Value = "'9.26"
$ws.Range("E21").Value = "  -1.79%  "
$ws.Range("D22").Value = "'439.83"
$ws.Range("E22").Value = "  +1.40%  "
$ws.Range("E23").Value = "  +6.73%  "
$ws.Range("D24").Value = "'7.14"
$ws.Range("E24").Value = "  +4.68%  "
$ws.Range("D25").Value = "'5.28"
$ws.Range("E25").Value = "  +10.68%  "
$ws.Range("E26").Value = "  +1.59%  "
$ws.Range("D27").Value = "3.344.56"
$ws.Range("E27").Value = "  +0.48%  "
$ws.Range("E28").Value = "  +0.27%  "
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  -0.37%  "
$ws.Range("E30").Value = "  +10.08%  "
$ws.Range("E31").Value = "  +2.76%  "
$ws.Range("E32").Value = "  +0.70%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "'1.53"
$ws.Range("E33").Value = "  +2.56%  "
$ws.Range("B34").Value = "Bittensor"
$ws.Range("C34").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D34").Value = "'560.65"
$ws.Range("E34").Value = "  +7.49%  "
$ws.Range("E35").Value = "  +3.10%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "'0.153"
$ws.Range("E36").Value = "  +12.18%  "
$ws.Range("B37").Value = "Cronos"
$ws.Range("C37").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D37").Value = "'0.141"
$ws.Range("E37").Value = "  +30.73%  "
$ws.Range("D38").Value = "'23.14"
$ws.Range("E38").Value = "  +2.63%  "
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("D40").Value = "'0.416"
$ws.Range("E40").Value = "  +4.62%  "
$ws.Range("D41").Value = "'3.13"
$ws.Range("E41").Value = "  +21.58%  "
$ws.Range("D42").Value = "'6.02"
$ws.Range("E42").Value = "  +10.95%  "
$ws.Range("D43").Value = "'2.05"
$ws.Range("E43").Value = "  +14.87%  "
$ws.Range("D44").Value = "'20.75"
$ws.Range("E44").Value = "  +3.78%  "
$ws.Range("D45").Value = "'160.52"
$ws.Range("E45").Value = "  -2.33%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").Value = "'190.08"
$ws.Range("E47").Value = "  -3.69%  "
$ws.Range("E48").Value = "  +3.86%  "
$ws.Range("D49").Value = "'44.34"
$ws.Range("E49").Value = "  +3.14%  "
$ws.Range("D50").Value = "'0.786"
$ws.Range("E50").Value = "  -1.77%  "
$ws.Range("D51").Value = "'4.29"
$ws.Range("E51").Value = "  +4.43%  "
